$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list values to match the latest scrape.
# D-column (Price) cells may look numeric (e.g. "3.268"); Excel would
# otherwise auto-convert them to actual numbers. Force them to stay as
# plain text (matching the original inlineStr cells), then clear the
# temporary "@" text format so the cell style index is left untouched.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "25.922.95"
$ws.Range("E2").Value = "  -0.52%  "
Set-TextValue "D3" "1.643.13"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.51%  "
Set-TextValue "D5" "215.45"
$ws.Range("E5").Value = "  -0.09%  "
Set-TextValue "D6" "0.5053"
$ws.Range("E6").Value = "  +0.06%  "
Set-TextValue "D7" "1.005"
Set-TextValue "D8" "0.2569"
$ws.Range("E8").Value = "  -0.57%  "
Set-TextValue "D9" "0.06395"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  +0.81%  "
Set-TextValue "D11" "0.07794"
$ws.Range("E11").Value = "  +0.66%  "
Set-TextValue "D12" "1.659.24"
$ws.Range("E12").Value = "  +0.79%  "
Set-TextValue "D13" "4.279"
$ws.Range("E13").Value = "  +0.61%  "
Set-TextValue "D14" "1.868.76"
$ws.Range("E14").Value = "  -0.13%  "
Set-TextValue "D15" "0.5431"
$ws.Range("E15").Value = "  -0.49%  "
Set-TextValue "D16" "0.0₅7860"
$ws.Range("E16").Value = "  -0.82%  "
Set-TextValue "D17" "64.84"
Set-TextValue "D18" "25.962.18"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E20").Value = "  -2.56%  "
Set-TextValue "D21" "4.392"
$ws.Range("E21").Value = "  +1.91%  "
Set-TextValue "D22" "9.977"
$ws.Range("E22").Value = "  -0.36%  "
Set-TextValue "D23" "5.981"
$ws.Range("E23").Value = "  +0.27%  "
Set-TextValue "D24" "1.006"
$ws.Range("E24").Value = "  -0.55%  "
Set-TextValue "D25" "1.870"
$ws.Range("E25").Value = "  -3.23%  "
Set-TextValue "D26" "139.93"
$ws.Range("E26").Value = "  -1.36%  "
Set-TextValue "D27" "0.1144"
$ws.Range("E27").Value = "  -1.21%  "
Set-TextValue "D28" "6.854"
$ws.Range("E28").Value = "  +1.49%  "
Set-TextValue "D29" "15.73"
$ws.Range("E29").Value = "  -0.06%  "
Set-TextValue "D30" "1.244"
$ws.Range("E30").Value = "  +0.10%  "
Set-TextValue "D31" "0.04924"
$ws.Range("E31").Value = "  -2.82%  "
Set-TextValue "D32" "3.268"
$ws.Range("E32").Value = "  +0.19%  "
Set-TextValue "D33" "3.198"
$ws.Range("E33").Value = "  +0.04%  "
Set-TextValue "D34" "1.535"
$ws.Range("E34").Value = "  -0.60%  "
Set-TextValue "D35" "2.372"
$ws.Range("E35").Value = "  +1.14%  "
Set-TextValue "D36" "0.8945"
$ws.Range("E36").Value = "  -0.23%  "
Set-TextValue "D37" "2.608"
$ws.Range("E37").Value = "  -0.50%  "
Set-TextValue "D38" "1.143.54"
$ws.Range("E38").Value = "  -0.39%  "
Set-TextValue "D39" "0.5554"
$ws.Range("E39").Value = "  -1.51%  "
Set-TextValue "D40" "0.01561"
$ws.Range("E40").Value = "  -0.78%  "
Set-TextValue "D41" "1.006"
$ws.Range("E41").Value = "  -0.53%  "
Set-TextValue "D42" "5.693"
$ws.Range("E42").Value = "  +0.30%  "
Set-TextValue "D43" "0.8174"
$ws.Range("E43").Value = "  +0.17%  "
Set-TextValue "D44" "99.40"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D45" "1.779.05"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D46" "0.0₈119"
$ws.Range("E46").Value = "  +4.76%  "
Set-TextValue "D47" "0.4527"
$ws.Range("E47").Value = "  -0.27%  "
Set-TextValue "D48" "55.38"
$ws.Range("E48").Value = "  +0.63%  "
Set-TextValue "D49" "1.006"
$ws.Range("E49").Value = "  -0.64%  "
Set-TextValue "D50" "0.05079"
$ws.Range("E50").Value = "  +0.69%  "
Set-TextValue "D51" "1.006"
$ws.Range("E51").Value = "  -0.45%  "
